$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from an existing header cell (H1) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-18
$values = @{
    2  = @(8, 8)
    3  = @(8, 8)
    4  = @(6, 6)
    5  = @(4, 4)
    6  = @(3, 5)
    7  = @(8, 8)
    8  = @(9, 9)
    9  = @(9, 9)
    10 = @(4, 5)
    11 = @(7, 7)
    12 = @(4, 4)
    13 = @(5, 5)
    14 = @(5, 5)
    15 = @(2, 3)
    16 = @(4, 5)
    17 = @(5, 5)
    18 = @(5, 5)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

$wb.Save()
